# The diff updates the classification-report style metrics tables on the
# "SGD", "LinearSVC", "MLP Neural Network" and "Gaussian Process" sheets.
# Support counts (row 5, columns B/C/D) changed from 7/22/9 to 12/15/11
# (total still 38) and the precision/recall/f1-score/accuracy values were
# recomputed accordingly for each model.

$wb = $excel.ActiveWorkbook

function Set-Metrics {
    param($SheetName, $Values)

    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($cellRef in $Values.Keys) {
        $ws.Range($cellRef).Value = $Values[$cellRef]
    }
}

# --- SGD ---
Set-Metrics "SGD" @{
    "C2" = 0.9333333333333333
    "D2" = 0.9090909090909091
    "E2" = 0.9473684210526315
    "F2" = 0.9474747474747475
    "G2" = 0.9473684210526315

    "C3" = 0.9333333333333333
    "D3" = 0.9090909090909091
    "E3" = 0.9473684210526315
    "F3" = 0.9474747474747475
    "G3" = 0.9473684210526315

    "C4" = 0.9333333333333333
    "D4" = 0.9090909090909091
    "E4" = 0.9473684210526315
    "F4" = 0.9474747474747475
    "G4" = 0.9473684210526315

    "B5" = 12
    "C5" = 15
    "D5" = 11
    "E5" = 0.9473684210526315
}

# --- LinearSVC ---
Set-Metrics "LinearSVC" @{
    "C2" = 0.9333333333333333
    "D2" = 0.9090909090909091
    "E2" = 0.9473684210526315
    "F2" = 0.9474747474747475
    "G2" = 0.9473684210526315

    "C3" = 0.9333333333333333
    "D3" = 0.9090909090909091
    "E3" = 0.9473684210526315
    "F3" = 0.9474747474747475
    "G3" = 0.9473684210526315

    "C4" = 0.9333333333333333
    "D4" = 0.9090909090909091
    "E4" = 0.9473684210526315
    "F4" = 0.9474747474747475
    "G4" = 0.9473684210526315

    "B5" = 12
    "C5" = 15
    "D5" = 11
    "E5" = 0.9473684210526315
}

# --- MLP Neural Network ---
Set-Metrics "MLP Neural Network" @{
    "D2" = 0.9166666666666666
    "E2" = 0.9736842105263158
    "F2" = 0.9722222222222222
    "G2" = 0.975877192982456

    "C3" = 0.9333333333333333
    "E3" = 0.9736842105263158
    "F3" = 0.9777777777777779
    "G3" = 0.9736842105263158

    "C4" = 0.9655172413793104
    "D4" = 0.9565217391304348
    "E4" = 0.9736842105263158
    "F4" = 0.9740129935032483
    "G4" = 0.9738025723980116

    "B5" = 12
    "C5" = 15
    "D5" = 11
    "E5" = 0.9736842105263158
}

# --- Gaussian Process ---
Set-Metrics "Gaussian Process" @{
    "C2" = 0
    "D2" = 0.4230769230769231
    "E2" = 0.6052631578947368
    "F2" = 0.4743589743589744
    "G2" = 0.4382591093117409

    "C3" = 0
    "E3" = 0.6052631578947368
    "F3" = 0.6666666666666666
    "G3" = 0.6052631578947368

    "C4" = 0
    "D4" = 0.5945945945945945
    "E4" = 0.6052631578947368
    "F4" = 0.5315315315315315
    "G4" = 0.4879089615931721

    "B5" = 12
    "C5" = 15
    "D5" = 11
    "E5" = 0.6052631578947368
}
